$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section: request_type table, appended after the existing failure_code table
# Header row (row 10)
$ws.Range("A10").Value = "request_type"
$ws.Range("B10").Value = "值"
$ws.Range("C10").Value = "含义"

# Data row (row 11) - REQ_TYPE_BASIC_INFO
$ws.Range("A11").Value = "REQ_TYPE_BASIC_INFO"
$ws.Range("B11").Value = "0x00"
$ws.Range("C11").Value = "获取当前设备基本信息"

# Data row (row 12) - REQ_TYPE_AP_LIST
$ws.Range("A12").Value = "REQ_TYPE_AP_LIST"
$ws.Range("B12").Value = "0x01"
$ws.Range("C12").Value = "获取当前设备检测到的AP列表"

# Update selection to match target state (new block A10:C12, active cell C12)
$ws.Range("A10:C12").Select()
